# Apply scheduled market-data refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 47084
$ws.Range("J126").Value = 47084
$ws.Range("L126").Value = 47084
$ws.Range("N126").Value = -56964
$ws.Range("H128").Value = 44041.168
$ws.Range("J128").Value = 44041.168
$ws.Range("L128").Value = 44041.168
$ws.Range("N128").Value = -54001.168
$ws.Range("H130").Value = 44803.2
$ws.Range("J130").Value = 44803.2
$ws.Range("L130").Value = 44803.2
$ws.Range("N130").Value = -54843.2
$ws.Range("H135").Value = 16130410
$ws.Range("I135").Value = 1383.4445
$ws.Range("J135").Value = 38462908
$ws.Range("K135").Value = 12451.0005
$ws.Range("L135").Value = 346166172
$ws.Range("M135").Value = -9916.0005
$ws.Range("N135").Value = -346171242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 49803
$ws.Range("J118").Value = 49803
$ws.Range("L118").Value = 49803
$ws.Range("N118").Value = -53117
$ws.Range("H123").Value = 35610.5
$ws.Range("J123").Value = 35610.5
$ws.Range("L123").Value = 35610.5
$ws.Range("N123").Value = -45410.5
$ws.Range("H125").Value = 48930.668
$ws.Range("J125").Value = 48930.668
$ws.Range("L125").Value = 48930.668
$ws.Range("N125").Value = -58770.668
$ws.Range("H138").Value = 41788.75
$ws.Range("J138").Value = 41788.75
$ws.Range("L138").Value = 41788.75
$ws.Range("N138").Value = -52068.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2667
$ws.Range("I86").Value = 2659.8
$ws.Range("J86").Value = 2681.4
$ws.Range("K86").Value = 2659.8
$ws.Range("L86").Value = 2681.4
$ws.Range("M86").Value = -1536.8
$ws.Range("N86").Value = -4927.4
$ws.Range("H89").Value = 2667
$ws.Range("I89").Value = 2659.8
$ws.Range("J89").Value = 2681.4
$ws.Range("K89").Value = 13299
$ws.Range("L89").Value = 13407
$ws.Range("M89").Value = -7683
$ws.Range("N89").Value = -24639
$ws.Range("H94").Value = 1001.3461
$ws.Range("I94").Value = 890.3333
$ws.Range("J94").Value = 1467.6
$ws.Range("K94").Value = 890.3333
$ws.Range("L94").Value = 1467.6
$ws.Range("M94").Value = -439.3333
$ws.Range("N94").Value = -2369.6
$ws.Range("H124").Value = 47992
$ws.Range("J124").Value = 47992
$ws.Range("L124").Value = 47992
$ws.Range("N124").Value = -57812
$ws.Range("H125").Value = 50172
$ws.Range("J125").Value = 50172
$ws.Range("L125").Value = 50172
$ws.Range("N125").Value = -60012
$ws.Range("H126").Value = 41749
$ws.Range("J126").Value = 41749
$ws.Range("L126").Value = 41749
$ws.Range("N126").Value = -51629
$ws.Range("H130").Value = 47115
$ws.Range("J130").Value = 47115
$ws.Range("L130").Value = 47115
$ws.Range("N130").Value = -57155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49003.75
$ws.Range("J20").Value = 49003.75
$ws.Range("L20").Value = 49003.75
$ws.Range("N20").Value = -49475.75
$ws.Range("H30").Value = 49003.75
$ws.Range("J30").Value = 49003.75
$ws.Range("L30").Value = 49003.75
$ws.Range("N30").Value = -49185.75
$ws.Range("H58").Value = 1748.1613
$ws.Range("I58").Value = 1535.5714
$ws.Range("J58").Value = 3732.3333
$ws.Range("K58").Value = 1535.5714
$ws.Range("L58").Value = 3732.3333
$ws.Range("M58").Value = -1332.5714
$ws.Range("N58").Value = -4138.3333
$ws.Range("H100").Value = 38174
$ws.Range("J100").Value = 47761
$ws.Range("L100").Value = 47761
$ws.Range("N100").Value = -49925
$ws.Range("H128").Value = 49003.75
$ws.Range("J128").Value = 49003.75
$ws.Range("L128").Value = 49003.75
$ws.Range("N128").Value = -58963.75
$ws.Range("H136").Value = 1748.1613
$ws.Range("I136").Value = 1535.5714
$ws.Range("J136").Value = 3732.3333
$ws.Range("K136").Value = 4606.7142
$ws.Range("L136").Value = 11196.9999
$ws.Range("M136").Value = -2056.7142
$ws.Range("N136").Value = -16296.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2400
$ws.Range("J104").Value = 2400
$ws.Range("L104").Value = 7200
$ws.Range("N104").Value = -12442

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 24574.666
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 24574.666
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 24574.666
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -26214.666
$ws.Range("H110").Value = 47702
$ws.Range("J110").Value = 47702
$ws.Range("L110").Value = 47702
$ws.Range("N110").Value = -55882
$ws.Range("H130").Value = 45784
$ws.Range("J130").Value = 45784
$ws.Range("L130").Value = 45784
$ws.Range("N130").Value = -55824
$ws.Range("H138").Value = 42872.727
$ws.Range("J138").Value = 42872.727
$ws.Range("L138").Value = 42872.727
$ws.Range("N138").Value = -53152.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5953473
$ws.Range("I82").Value = 1107.25
$ws.Range("J82").Value = 41667668
$ws.Range("K82").Value = 1107.25
$ws.Range("L82").Value = 41667668
$ws.Range("M82").Value = -746.25
$ws.Range("N82").Value = -41668390
$ws.Range("H85").Value = 5953473
$ws.Range("I85").Value = 1107.25
$ws.Range("J85").Value = 41667668
$ws.Range("K85").Value = 1107.25
$ws.Range("L85").Value = 41667668
$ws.Range("M85").Value = 140.75
$ws.Range("N85").Value = -41670164
$ws.Range("H127").Value = 50627.668
$ws.Range("J127").Value = 50627.668
$ws.Range("L127").Value = 50627.668
$ws.Range("N127").Value = -60547.668
$ws.Range("H130").Value = 37996
$ws.Range("J130").Value = 37996
$ws.Range("L130").Value = 37996
$ws.Range("N130").Value = -48036

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 100000000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H120").Value = 35206
$ws.Range("J120").Value = 35206
$ws.Range("L120").Value = 35206
$ws.Range("N120").Value = -44882
$ws.Range("H122").Value = 2041951.4
$ws.Range("I122").Value = 3175846.8
$ws.Range("K122").Value = 9527540.399999999
$ws.Range("M122").Value = -9525090.399999999
$ws.Range("H124").Value = 31857.25
$ws.Range("J124").Value = 31857.25
$ws.Range("L124").Value = 31857.25
$ws.Range("N124").Value = -41677.25
$ws.Range("H126").Value = 1133007.5
$ws.Range("I126").Value = 1338372.2
$ws.Range("K126").Value = 4015116.6
$ws.Range("M126").Value = -4012646.6
$ws.Range("H128").Value = 49715
$ws.Range("J128").Value = 49715
$ws.Range("L128").Value = 49715
$ws.Range("N128").Value = -59675
$ws.Range("H136").Value = 22781.318
$ws.Range("I136").Value = 50844.8
$ws.Range("J136").Value = 1993.5555
$ws.Range("K136").Value = 152534.4
$ws.Range("L136").Value = 5980.666499999999
$ws.Range("M136").Value = -149984.4
$ws.Range("N136").Value = -11080.6665
$ws.Range("H137").Value = 50999.5
$ws.Range("J137").Value = 50999.5
$ws.Range("L137").Value = 50999.5
$ws.Range("N137").Value = -61199.5
